# Atualizado por script em 01-12-2023 20:45
#
# The source scraper re-ran and the match listing shifted: a few adjacent
# fixtures swapped positions in the sheet and one new fixture (Al Nasr vs
# Al Wahda) was appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 2 <-> 3 swap (Al Wahda-Al Bataeh <-> Ittihad Kalba-Al Sharjah)
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "Ittihad Kalba"
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = "Al Sharjah"
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 3.85
$ws.Range("L2").Value = 4.06
$ws.Range("M2").Value = "18/08/2023 15:50"
$ws.Range("N2").Value = 3.67
$ws.Range("P2").Value = 3.72
$ws.Range("Q2").Value = "18/08/2023 15:50"
$ws.Range("R2").Value = 1.91
$ws.Range("T2").Value = 1.88
$ws.Range("U2").Value = "18/08/2023 15:50"
$ws.Range("V2").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ittihad-kalba-al-sharjah/l8ZYUBd2/"

$ws.Range("F3").Value = "Al Wahda"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "Al Bataeh"
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 1.22
$ws.Range("L3").Value = 1.24
$ws.Range("M3").Value = "18/08/2023 15:53"
$ws.Range("N3").Value = 6.71
$ws.Range("P3").Value = 6.55
$ws.Range("Q3").Value = "18/08/2023 15:58"
$ws.Range("R3").Value = 10.86
$ws.Range("T3").Value = 10.1
$ws.Range("U3").Value = "18/08/2023 15:58"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-al-bataeh/063NyrVE/"

# ---------------------------------------------------------------------
# Rows 18 <-> 19 swap (Al Ain-Ajman <-> Khorfakkan-Emirates Club)
# ---------------------------------------------------------------------
$ws.Range("F18").Value = "Khorfakkan"
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = "Emirates Club"
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 1.76
$ws.Range("K18").Value = "19/09/2023 16:42"
$ws.Range("L18").Value = 2.07
$ws.Range("M18").Value = "23/09/2023 17:54"
$ws.Range("N18").Value = 4.01
$ws.Range("O18").Value = "19/09/2023 16:42"
$ws.Range("P18").Value = 3.93
$ws.Range("Q18").Value = "23/09/2023 17:54"
$ws.Range("R18").Value = 4.14
$ws.Range("S18").Value = "19/09/2023 16:42"
$ws.Range("T18").Value = 3.26
$ws.Range("U18").Value = "23/09/2023 17:54"
$ws.Range("V18").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/khorfakkan-emirates-club/pYypuHy7/"

$ws.Range("F19").Value = "Al Ain"
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = "Ajman"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1.25
$ws.Range("K19").Value = "16/09/2023 18:13"
$ws.Range("L19").Value = 1.24
$ws.Range("M19").Value = "23/09/2023 17:52"
$ws.Range("N19").Value = 6.4
$ws.Range("O19").Value = "16/09/2023 18:13"
$ws.Range("P19").Value = 7
$ws.Range("Q19").Value = "23/09/2023 17:52"
$ws.Range("R19").Value = 9.51
$ws.Range("S19").Value = "16/09/2023 18:13"
$ws.Range("T19").Value = 9.550000000000001
$ws.Range("U19").Value = "23/09/2023 17:52"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-ain-ajman/0tKW2KTD/"

# ---------------------------------------------------------------------
# Rows 58 <-> 59 swap (Ajman-Hatta <-> Emirates Club-Ittihad Kalba)
# ---------------------------------------------------------------------
$ws.Range("F58").Value = "Emirates Club"
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = "Ittihad Kalba"
$ws.Range("I58").Value = 4
$ws.Range("J58").Value = 3.65
$ws.Range("L58").Value = 3.74
$ws.Range("M58").Value = "01/12/2023 13:41"
$ws.Range("N58").Value = 3.79
$ws.Range("P58").Value = 4.3
$ws.Range("Q58").Value = "01/12/2023 13:44"
$ws.Range("R58").Value = 1.86
$ws.Range("T58").Value = 1.83
$ws.Range("U58").Value = "01/12/2023 13:43"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/emirates-club-ittihad-kalba/ncfNFf9l/"

$ws.Range("F59").Value = "Ajman"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = "Hatta"
$ws.Range("I59").Value = 1
$ws.Range("J59").Value = 1.67
$ws.Range("L59").Value = 1.86
$ws.Range("M59").Value = "01/12/2023 13:42"
$ws.Range("N59").Value = 3.9
$ws.Range("P59").Value = 4.23
$ws.Range("Q59").Value = "01/12/2023 13:42"
$ws.Range("R59").Value = 4.5
$ws.Range("T59").Value = 3.65
$ws.Range("U59").Value = "01/12/2023 13:42"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/ajman-hatta/86eJGzfr/"

# ---------------------------------------------------------------------
# New row 60: Al Nasr vs Al Wahda (appended fixture)
# Copy row 59's formatting first (bold/bordered index column style,
# date-time number format on the match-date column) then fill in values.
# ---------------------------------------------------------------------
$ws.Range("A59:V59").Copy()
$ws.Range("A60:V60").PasteSpecial(-4122)

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "united-arab-emirates"
$ws.Range("C60").Value = "uae-league"
$ws.Range("D60").Value = "2023-2024"
$ws.Range("E60").Value = 45261.6875
$ws.Range("F60").Value = "Al Nasr"
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = "Al Wahda"
$ws.Range("I60").Value = 1
$ws.Range("J60").Value = 3.41
$ws.Range("K60").Value = "26/11/2023 14:13"
$ws.Range("L60").Value = 3.43
$ws.Range("M60").Value = "01/12/2023 16:25"
$ws.Range("N60").Value = 3.74
$ws.Range("O60").Value = "26/11/2023 14:13"
$ws.Range("P60").Value = 3.71
$ws.Range("Q60").Value = "01/12/2023 16:22"
$ws.Range("R60").Value = 2.02
$ws.Range("S60").Value = "26/11/2023 14:13"
$ws.Range("T60").Value = 2.07
$ws.Range("U60").Value = "01/12/2023 16:25"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-nasr-al-wahda/E5woZrA9/"
